$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update projection outputs per row (Year 1-5), fixing termination
# compensation proration logic. Values below reflect the corrected
# calculation results for headcount, eligibility, participation,
# deferral rates, contributions, compensation, and cost percentages.

# Row 2
$ws.Range("C2").Value = 102
$ws.Range("E2").Value = 0.8529411764705882
$ws.Range("G2").Value = 0.09823822539669926
$ws.Range("H2").Value = 0.08379142754424349
$ws.Range("I2").Value = 461915.3705097084
$ws.Range("J2").Value = 168343.6854558542
$ws.Range("L2").Value = 168343.6854558542
$ws.Range("M2").Value = 630259.0559655628
$ws.Range("N2").Value = 10165519.7688
$ws.Range("O2").Value = 9757778.838699998
$ws.Range("P2").Value = 0.01656026344786958
$ws.Range("Q2").Value = 0.01725225466150062

# Row 3
$ws.Range("B3").Value = 106
$ws.Range("C3").Value = 106
$ws.Range("D3").Value = 91
$ws.Range("E3").Value = 0.8584905660377359
$ws.Range("F3").Value = 0.8584905660377359
$ws.Range("G3").Value = 0.09607362778140782
$ws.Range("H3").Value = 0.08247830309535956
$ws.Range("I3").Value = 493240.0544632261
$ws.Range("J3").Value = 181407.607438643
$ws.Range("L3").Value = 181407.607438643
$ws.Range("M3").Value = 674647.6619018689
$ws.Range("N3").Value = 10801469.936564
$ws.Range("O3").Value = 10394096.778561
$ws.Range("P3").Value = 0.01679471484011274
$ws.Range("Q3").Value = 0.01745294577329862

# Row 4
$ws.Range("B4").Value = 109
$ws.Range("C4").Value = 109
$ws.Range("D4").Value = 94
$ws.Range("E4").Value = 0.8623853211009175
$ws.Range("F4").Value = 0.8623853211009175
$ws.Range("G4").Value = 0.09422337361318721
$ws.Range("H4").Value = 0.08125685430862017
$ws.Range("I4").Value = 526514.418027284
$ws.Range("J4").Value = 190958.4166268829
$ws.Range("L4").Value = 190958.4166268829
$ws.Range("M4").Value = 717472.834654167
$ws.Range("N4").Value = 11345286.94306092
$ws.Range("O4").Value = 10936942.59031783
$ws.Range("P4").Value = 0.01683151934237134
$ws.Range("Q4").Value = 0.017459945048622

# Row 5
$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 110
$ws.Range("D5").Value = 95
$ws.Range("E5").Value = 0.8636363636363636
$ws.Range("F5").Value = 0.8636363636363636
$ws.Range("G5").Value = 0.09389805491180749
$ws.Range("H5").Value = 0.08109377469656102
$ws.Range("I5").Value = 548242.7938357895
$ws.Range("J5").Value = 199248.6407595329
$ws.Range("L5").Value = 199248.6407595329
$ws.Range("M5").Value = 747491.4345953225
$ws.Range("N5").Value = 11786101.98315275
$ws.Range("O5").Value = 11375407.29982737
$ws.Range("P5").Value = 0.016905389164657
$ws.Range("Q5").Value = 0.01751573684421451

# Row 6
$ws.Range("B6").Value = 112
$ws.Range("C6").Value = 112
$ws.Range("D6").Value = 97
$ws.Range("E6").Value = 0.8660714285714286
$ws.Range("F6").Value = 0.8660714285714286
$ws.Range("G6").Value = 0.09400854608248253
$ws.Range("H6").Value = 0.08141811580357862
$ws.Range("I6").Value = 574500.8661924924
$ws.Range("J6").Value = 209316.4942531334
$ws.Range("L6").Value = 209316.4942531334
$ws.Range("M6").Value = 783817.3604456257
$ws.Range("N6").Value = 12106775.67824733
$ws.Range("O6").Value = 11692310.15442219
$ws.Range("P6").Value = 0.01728920232900819
$ws.Range("Q6").Value = 0.01790206481770133
